$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9320.61
$ws.Range("B10").Value = 9390.1
$ws.Range("C10").Value = 286
$ws.Range("D10").Value = 283.87
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = -0.74
$ws.Range("G10").Value = 42612.673043981478
$ws.Range("G10").NumberFormat = "m/d/yy h:mm"
$ws.Range("H10").Value = $false
